$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column (C) for rows 2-5
# from 2023-11-13 (serial 45243) to 2023-11-14 (serial 45244)
$ws.Range("C2:C5").Value = 45244
